$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Update time_taken column (F2:F13) on the "data" sheet ---
$data.Range("F2").Value = "2021-10-05 14:19:51.204816"
$data.Range("F3").Value = "2021-10-05 14:19:51.204824"
$data.Range("F4").Value = "2021-10-05 14:19:51.204827"
$data.Range("F5").Value = "2021-10-05 14:19:51.204830"
$data.Range("F6").Value = "2021-10-05 14:19:51.204833"
$data.Range("F7").Value = "2021-10-05 14:19:51.204836"
$data.Range("F8").Value = "2021-10-05 14:19:51.204838"
$data.Range("F9").Value = "2021-10-05 14:19:51.204841"
$data.Range("F10").Value = "2021-10-05 14:19:51.204844"
$data.Range("F11").Value = "2021-10-05 14:19:51.204846"
$data.Range("F12").Value = "2021-10-05 14:19:51.204849"
$data.Range("F13").Value = "2021-10-05 14:19:51.204851"

# --- Add the new "metadata" sheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Reuse the "data" sheet's bold/bordered/centered header style (style index 1)
# for the metadata header row, via copy/paste-format so no duplicate style is
# created in the styles table.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Row index cell A2 reuses the same style as data!A2
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Cutaneous photosensitivity with a likely genetic cause"
$meta.Range("C2").Value = 560

# "1.8" needs to stay plain text (not get reinterpreted as the number 1.8).
# Stage it in a scratch cell formatted as Text, copy only the VALUE across
# (xlPasteValues) so the destination cell keeps the default (no) style, then
# clean up the scratch cell.
$meta.Range("Z1").NumberFormat = "@"
$meta.Range("Z1").Value = "1.8"
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

$meta.Range("E2").Value = "2021-07-12T13:21:34.496460Z"
$meta.Range("F2").Value = "2021-10-05 14:19:51.201060"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/560/?format=json"

# Keep "data" as the active tab (the diff leaves <bookViews>/activeTab
# untouched, i.e. still pointing at the first/original sheet).
$data.Activate()
